# Applies the diff for Jogos_da_Semana_FlashScore_2025-01-28.xlsx (Sheet1)
# Updates odds values in rows 5, 8, 10-12, 15, 17-19, 21-27, 32-33

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5
$ws.Range("G5").Value = 2.38
$ws.Range("I5").Value = 3.5
$ws.Range("J5").Value = 3.25
$ws.Range("L5").Value = 4.33
$ws.Range("Y5").Value = 5.5
$ws.Range("Z5").Value = 9.5
$ws.Range("AA5").Value = 10
$ws.Range("AB5").Value = 23
$ws.Range("AC5").Value = 23
$ws.Range("AF5").Value = 6
$ws.Range("AJ5").Value = 7.5
$ws.Range("AK5").Value = 15
$ws.Range("AL5").Value = 13
$ws.Range("AM5").Value = 41
$ws.Range("AO5").Value = 51

# Row 8
$ws.Range("G8").Value = 1.85
$ws.Range("H8").Value = 3.7
$ws.Range("I8").Value = 3.3
$ws.Range("J8").Value = 2.5
$ws.Range("L8").Value = 4
$ws.Range("U8").Value = 1.3
$ws.Range("V8").Value = 3.4
$ws.Range("W8").Value = 1.62
$ws.Range("X8").Value = 2.2
$ws.Range("Y8").Value = 9.5
$ws.Range("Z8").Value = 11
$ws.Range("AA8").Value = 8.5
$ws.Range("AB8").Value = 17
$ws.Range("AC8").Value = 15
$ws.Range("AD8").Value = 21
$ws.Range("AF8").Value = 7.5
$ws.Range("AG8").Value = 13
$ws.Range("AI8").Value = 151
$ws.Range("AJ8").Value = 13
$ws.Range("AK8").Value = 21
$ws.Range("AL8").Value = 12
$ws.Range("AM8").Value = 41
$ws.Range("AN8").Value = 26
$ws.Range("AO8").Value = 29

# Row 10
$ws.Range("G10").Value = 3.1
$ws.Range("H10").Value = 3.5
$ws.Range("I10").Value = 2.05
$ws.Range("J10").Value = 3.75
$ws.Range("L10").Value = 2.75
$ws.Range("Y10").Value = 12
$ws.Range("Z10").Value = 17
$ws.Range("AA10").Value = 12
$ws.Range("AB10").Value = 34
$ws.Range("AC10").Value = 23
$ws.Range("AD10").Value = 29
$ws.Range("AG10").Value = 13
$ws.Range("AJ10").Value = 9
$ws.Range("AK10").Value = 11
$ws.Range("AL10").Value = 9
$ws.Range("AM10").Value = 19
$ws.Range("AN10").Value = 15

# Row 11
$ws.Range("G11").Value = 1.65
$ws.Range("H11").Value = 3.6
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("AF11").Value = 7
$ws.Range("AM11").Value = 51
$ws.Range("AP11").Value = 1.65
$ws.Range("AQ11").Value = 2.2

# Row 12
$ws.Range("G12").Value = 2.32
$ws.Range("H12").Value = 2.7
$ws.Range("I12").Value = 3.55
$ws.Range("J12").Value = 3.05
$ws.Range("K12").Value = 1.85
$ws.Range("L12").Value = 4.2
$ws.Range("M12").Value = 1.13
$ws.Range("N12").Value = 5
$ws.Range("O12").Value = 1.55
$ws.Range("P12").Value = 2.3
$ws.Range("Q12").Value = 2.6
$ws.Range("R12").Value = 1.44
$ws.Range("S12").Value = 4.7
$ws.Range("T12").Value = 1.15
$ws.Range("U12").Value = 1.57
$ws.Range("V12").Value = 2.25
$ws.Range("W12").Value = 2.1
$ws.Range("X12").Value = 1.65
$ws.Range("Y12").Value = 5.6
$ws.Range("AB12").Value = 25
$ws.Range("AC12").Value = 24
$ws.Range("AD12").Value = 45
$ws.Range("AE12").Value = 5
$ws.Range("AF12").Value = 5.4
$ws.Range("AG12").Value = 17
$ws.Range("AH12").Value = 110
$ws.Range("AJ12").Value = 7.8
$ws.Range("AL12").Value = 12.5
$ws.Range("AN12").Value = 40
$ws.Range("AO12").Value = 55

# Row 15
$ws.Range("G15").Value = 1.38
$ws.Range("I15").Value = 7
$ws.Range("L15").Value = 7.5
$ws.Range("W15").Value = 2.05
$ws.Range("X15").Value = 1.7
$ws.Range("Y15").Value = 6.5
$ws.Range("AB15").Value = 9
$ws.Range("AD15").Value = 29
$ws.Range("AI15").Value = 451
$ws.Range("AK15").Value = 41
$ws.Range("AL15").Value = 21

# Row 17
$ws.Range("I17").Value = 6.25

# Row 18
$ws.Range("G18").Value = 3.3
$ws.Range("I18").Value = 2.15
$ws.Range("N18").Value = 9
$ws.Range("AL18").Value = 9.5
$ws.Range("AM18").Value = 21

# Row 19
$ws.Range("G19").Value = 2.45
$ws.Range("I19").Value = 2.88
$ws.Range("J19").Value = 3.2
$ws.Range("Y19").Value = 8
$ws.Range("Z19").Value = 12
$ws.Range("AA19").Value = 10
$ws.Range("AK19").Value = 13

# Row 21
$ws.Range("M21").Value = 1.05
$ws.Range("N21").Value = 11
$ws.Range("O21").Value = 1.29
$ws.Range("P21").Value = 3.5
$ws.Range("Q21").Value = 2
$ws.Range("R21").Value = 1.85
$ws.Range("S21").Value = 3.4
$ws.Range("T21").Value = 1.3

# Row 22
$ws.Range("G22").Value = 1.67
$ws.Range("M22").Value = 1.05
$ws.Range("N22").Value = 11
$ws.Range("O22").Value = 1.3
$ws.Range("P22").Value = 3.4
$ws.Range("Q22").Value = 2.03
$ws.Range("R22").Value = 1.83
$ws.Range("AH22").Value = 51

# Row 23
$ws.Range("J23").Value = 2.05
$ws.Range("K23").Value = 2.2
$ws.Range("L23").Value = 7.5
$ws.Range("M23").Value = 1.07
$ws.Range("N23").Value = 9
$ws.Range("O23").Value = 1.33
$ws.Range("P23").Value = 3.25
$ws.Range("Q23").Value = 2.05
$ws.Range("R23").Value = 1.75
$ws.Range("S23").Value = 3.75
$ws.Range("T23").Value = 1.25
$ws.Range("U23").Value = 1.44
$ws.Range("V23").Value = 2.63
$ws.Range("Y23").Value = 5.5
$ws.Range("AA23").Value = 9
$ws.Range("AE23").Value = 9
$ws.Range("AJ23").Value = 13

# Row 24
$ws.Range("H24").Value = 3.2
$ws.Range("I24").Value = 4.1
$ws.Range("L24").Value = 5
$ws.Range("M24").Value = 1.1
$ws.Range("N24").Value = 7
$ws.Range("Q24").Value = 2.6
$ws.Range("R24").Value = 1.48
$ws.Range("U24").Value = 1.57
$ws.Range("V24").Value = 2.25
$ws.Range("Y24").Value = 5.5
$ws.Range("Z24").Value = 8
$ws.Range("AE24").Value = 6.5
$ws.Range("AJ24").Value = 8
$ws.Range("AK24").Value = 19
$ws.Range("AL24").Value = 15
$ws.Range("AM24").Value = 41

# Row 25
$ws.Range("Q25").Value = 1.93
$ws.Range("R25").Value = 1.93

# Row 26
$ws.Range("H26").Value = 3
$ws.Range("I26").Value = 2.38
$ws.Range("J26").Value = 4
$ws.Range("K26").Value = 1.95
$ws.Range("N26").Value = 8
$ws.Range("O26").Value = 1.44
$ws.Range("P26").Value = 2.63
$ws.Range("AE26").Value = 7
$ws.Range("AH26").Value = 67
$ws.Range("AI26").Value = 451
$ws.Range("AK26").Value = 10
$ws.Range("AP26").Value = 1.78
$ws.Range("AQ26").Value = 2.1

# Row 27
$ws.Range("G27").Value = 1.8
$ws.Range("H27").Value = 3.5
$ws.Range("I27").Value = 4.33
$ws.Range("J27").Value = 2.5
$ws.Range("L27").Value = 4.75
$ws.Range("M27").Value = 1.06
$ws.Range("N27").Value = 10
$ws.Range("O27").Value = 1.3
$ws.Range("P27").Value = 3.4
$ws.Range("Q27").Value = 2
$ws.Range("R27").Value = 1.8
$ws.Range("S27").Value = 3.5
$ws.Range("T27").Value = 1.29
$ws.Range("W27").Value = 1.83
$ws.Range("X27").Value = 1.83
$ws.Range("Y27").Value = 7
$ws.Range("Z27").Value = 8.5
$ws.Range("AA27").Value = 8.5
$ws.Range("AB27").Value = 15
$ws.Range("AC27").Value = 15
$ws.Range("AD27").Value = 29
$ws.Range("AE27").Value = 9.5
$ws.Range("AG27").Value = 15
$ws.Range("AH27").Value = 51
$ws.Range("AI27").Value = 301
$ws.Range("AJ27").Value = 12
$ws.Range("AK27").Value = 21
$ws.Range("AL27").Value = 15
$ws.Range("AM27").Value = 41
$ws.Range("AN27").Value = 34
$ws.Range("AO27").Value = 41

# Row 32
$ws.Range("M32").Value = 1.1
$ws.Range("N32").Value = 7

# Row 33
$ws.Range("G33").Value = 2.63
$ws.Range("I33").Value = 2.55
$ws.Range("L33").Value = 3.2
$ws.Range("AK33").Value = 13
$ws.Range("AL33").Value = 11
$ws.Range("AM33").Value = 26
$ws.Range("AO33").Value = 34
